$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("G4").Value = 1493
$ws.Range("G7").Value = 24721
$ws.Range("K2").Value = 7149
$ws.Range("K3").Value = 7413
$ws.Range("K4").Value = 1545
$ws.Range("K6").Value = 8192
$ws.Range("K7").Value = 24823

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("G101").Value = 24721
$ws.Range("G63").Value = 296
$ws.Range("K10").Value = 144
$ws.Range("K101").Value = 24823
$ws.Range("K11").Value = 454
$ws.Range("K18").Value = 164
$ws.Range("K20").Value = 608
$ws.Range("K22").Value = 78
$ws.Range("K23").Value = 254
$ws.Range("K29").Value = 1365
$ws.Range("K33").Value = 1057
$ws.Range("K34").Value = 141
$ws.Range("K35").Value = 39
$ws.Range("K36").Value = 315
$ws.Range("K37").Value = 834
$ws.Range("K40").Value = 55
$ws.Range("K42").Value = 917
$ws.Range("K43").Value = 204
$ws.Range("K46").Value = 51
$ws.Range("K48").Value = 318
$ws.Range("K49").Value = 140
$ws.Range("K52").Value = 644
$ws.Range("K54").Value = 483
$ws.Range("K60").Value = 142
$ws.Range("K63").Value = 71
$ws.Range("K65").Value = 584
$ws.Range("K66").Value = 75
$ws.Range("K67").Value = 973
$ws.Range("K69").Value = 57
$ws.Range("K7").Value = 747
$ws.Range("K73").Value = 221
$ws.Range("K76").Value = 341
$ws.Range("K78").Value = 300
$ws.Range("K8").Value = 1616
$ws.Range("K80").Value = 93
$ws.Range("K83").Value = 527
$ws.Range("K84").Value = 196
$ws.Range("K85").Value = 1137
$ws.Range("K86").Value = 153
$ws.Range("K87").Value = 50
$ws.Range("K89").Value = 373
$ws.Range("K9").Value = 115
$ws.Range("K90").Value = 238
$ws.Range("K95").Value = 407
$ws.Range("K97").Value = 199
$ws.Range("K99").Value = 425

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 237
$ws.Range("K7").Value = 747

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 154
$ws.Range("K6").Value = 156
$ws.Range("K7").Value = 454

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K4").Value = 41
$ws.Range("K6").Value = 112
$ws.Range("K7").Value = 373

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 371
$ws.Range("K3").Value = 397
$ws.Range("K4").Value = 58
$ws.Range("K6").Value = 280
$ws.Range("K7").Value = 1137

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 174
$ws.Range("K3").Value = 177
$ws.Range("K6").Value = 238
$ws.Range("K7").Value = 644

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K3").Value = 489
$ws.Range("K4").Value = 93
$ws.Range("K7").Value = 1616

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 183
$ws.Range("K3").Value = 186
$ws.Range("K7").Value = 527

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 265
$ws.Range("K3").Value = 379
$ws.Range("K6").Value = 333
$ws.Range("K7").Value = 1057

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K3").Value = 140
$ws.Range("K6").Value = 92
$ws.Range("K7").Value = 407

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K6").Value = 249
$ws.Range("K7").Value = 834

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 190
$ws.Range("K3").Value = 138
$ws.Range("K4").Value = 24
$ws.Range("K7").Value = 584

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K2").Value = 111
$ws.Range("K3").Value = 178
$ws.Range("K6").Value = 104
$ws.Range("K7").Value = 425

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 271
$ws.Range("K3").Value = 349
$ws.Range("K6").Value = 274
$ws.Range("K7").Value = 973

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K3").Value = 82
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 196

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K2").Value = 75
$ws.Range("K6").Value = 266
$ws.Range("K7").Value = 483

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K3").Value = 488
$ws.Range("K7").Value = 1365

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K3").Value = 78
$ws.Range("K6").Value = 147
$ws.Range("K7").Value = 318

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 173
$ws.Range("K7").Value = 341

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 250
$ws.Range("K3").Value = 272
$ws.Range("K7").Value = 917

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K6").Value = 65
$ws.Range("K7").Value = 144

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K6").Value = 101
$ws.Range("K7").Value = 300

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K3").Value = 88
$ws.Range("K7").Value = 254

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 209
$ws.Range("K6").Value = 168
$ws.Range("K7").Value = 608

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 164

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 120
$ws.Range("K7").Value = 315

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("K2").Value = 55
$ws.Range("K7").Value = 141

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 39

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K3").Value = 40
$ws.Range("K7").Value = 115

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K3").Value = 57
$ws.Range("K6").Value = 75
$ws.Range("K7").Value = 221

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K3").Value = 45
$ws.Range("K6").Value = 104
$ws.Range("K7").Value = 199

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 153

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 86
$ws.Range("K7").Value = 238

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K2").Value = 48
$ws.Range("K7").Value = 142

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K2").Value = 41
$ws.Range("K7").Value = 204

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("K3").Value = 23
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 50

Write-Output "Applied 162 cell updates across 44 sheets for 2024-11-19 data"
